$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corregido error del step: update the "Componentes Conexos" values for
# each q step, and remove the extra rows (q > 1) that were left over
# from a wrong step size.

$ws.Range("B2").Value = 839255
$ws.Range("B3").Value = 631545
$ws.Range("B4").Value = 429537
$ws.Range("B5").Value = 245326
$ws.Range("B6").Value = 102841
$ws.Range("B7").Value = 35225
$ws.Range("B8").Value = 9972
$ws.Range("B9").Value = 1853
$ws.Range("B10").Value = 108

# Remove rows 12-26, leaving data only through row 11 (q = 1)
$ws.Range("A12:B26").EntireRow.Delete()
